# Actualización automática 2025-08-07 08:30:08
# Applies the updated sales figure (475.2) for ALMEIDA CUATIN JHONATHANN CARLOS
# / HERRERA CAICEDO LUIS FRANKLIN across the three report sheets, and the
# resulting roll-up totals / percentages.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---------------------------------------------
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("D15").Value = 475.2
$wsVentasGrupo.Range("D33").Value = "3 de 31"

# --- Sheet "VENTA MENSUAL" --------------------------------------------------
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F15").Value = 475.2
$wsVentaMensual.Range("F33").Value = 3926.34

# --- Sheet "CUMPLIMIENTO MENSUAL" ------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D3").Value = 1391.04
$wsCumplimiento.Range("E3").Value = 1729.0745
$wsCumplimiento.Range("F3").Value = 0.445829792464347

$wsCumplimiento.Range("D19").Value = 3926.34
$wsCumplimiento.Range("E19").Value = 28182.94107555788
$wsCumplimiento.Range("F19").Value = 0.1222805328702546
